$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# "BAL case split": row 37 (TCID 42, LeaveBalance_48EmployeeCreation) is split
# into 4 rows covering employee-creation sub-ranges 1-25 / 26-41 / 83-100 /
# 101-123, and every row below shifts down by three positions (renumbered).
# ---------------------------------------------------------------------------

# 1) Update the existing row 37 in place for the first split case (1_25).
$ws.Cells.Item(37, 4).Value2 = "com.darwinbox.leaves.Accural.Custom.LeaveBalance_48EmployeeCreation_1_25"
$ws.Cells.Item(37, 2).Value2 = "LeaveBalance_1_25"
$ws.Cells.Item(37, 3).Value2 = "LeaveBalance_1_25"

# 2) Insert three new rows (38:40) below it to hold the other split cases;
#    they inherit row 37's formatting (text-style on columns A/G).
$ws.Rows("38:40").Insert()

# Columns E (FileName), F (SheetName) and G (TestDataRow) are identical
# across all four split rows, same as the original row.
$ws.Cells.Item(38, 5).Value2 = "Accural//LeaveBalance.xlsx"
$ws.Cells.Item(38, 6).Value2 = "dummySheet"
$ws.Cells.Item(38, 7).Value2 = "All"

$ws.Cells.Item(39, 5).Value2 = "Accural//LeaveBalance.xlsx"
$ws.Cells.Item(39, 6).Value2 = "dummySheet"
$ws.Cells.Item(39, 7).Value2 = "All"

$ws.Cells.Item(40, 5).Value2 = "Accural//LeaveBalance.xlsx"
$ws.Cells.Item(40, 6).Value2 = "dummySheet"
$ws.Cells.Item(40, 7).Value2 = "All"

# Module / TestCaseDescription for the remaining three split cases.
$ws.Cells.Item(38, 2).Value2 = "LeaveBalance_26_41"
$ws.Cells.Item(38, 3).Value2 = "LeaveBalance_26_41"

$ws.Cells.Item(39, 2).Value2 = "LeaveBalance_83_100"
$ws.Cells.Item(39, 3).Value2 = "LeaveBalance_83_100"

$ws.Cells.Item(40, 2).Value2 = "LeaveBalance_101_123"
$ws.Cells.Item(40, 3).Value2 = "LeaveBalance_101_123"

# ClassName for the remaining three split cases.
$ws.Cells.Item(38, 4).Value2 = "com.darwinbox.leaves.Accural.Custom.LeaveBalance_48EmployeeCreation_26_41"
$ws.Cells.Item(39, 4).Value2 = "com.darwinbox.leaves.Accural.Custom.LeaveBalance_48EmployeeCreation_83_100"
$ws.Cells.Item(40, 4).Value2 = "com.darwinbox.leaves.Accural.Custom.LeaveBalance_48EmployeeCreation_101_123"

# 3) Renumber TCID for the split rows (42-45).
$ws.Cells.Item(37, 1).Value2 = "42"
$ws.Cells.Item(38, 1).Value2 = "43"
$ws.Cells.Item(39, 1).Value2 = "44"
$ws.Cells.Item(40, 1).Value2 = "45"

# 4) Renumber TCID for all the rows that shifted down by three (46-52).
$ws.Cells.Item(41, 1).Value2 = "46"
$ws.Cells.Item(42, 1).Value2 = "47"
$ws.Cells.Item(43, 1).Value2 = "48"
$ws.Cells.Item(44, 1).Value2 = "49"
$ws.Cells.Item(45, 1).Value2 = "50"
$ws.Cells.Item(46, 1).Value2 = "51"
$ws.Cells.Item(47, 1).Value2 = "52"

# 5) Restore the active selection to the last data row, as in the authored file.
$ws.Range("B47").Select()
